{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const paragraph = paragraphs.items[i];\n  // Tighten the bullet hanging indent: w:hanging=\"244\" (-12.2pt) -> w:hanging=\"210\" (-10.5pt)\n  paragraph.firstLineIndent = -10.5;\n}\nawait context.sync();\n\n// Drop the trailing space in the first bullet's text:\n// \"Unordered information. \" -> \"Unordered information.\"\nconst firstParagraph = paragraphs.items[0];\nfirstParagraph.load(\"text\");\nawait context.sync();\n\nconst trimmedText = firstParagraph.text.replace(/\\s+$/, \"\");\nif (trimmedText !== firstParagraph.text) {\n  firstParagraph.clear();\n  await context.sync();\n  firstParagraph.insertText(trimmedText, Word.InsertLocation.start);\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Tighten the bullet hanging indent on every list paragraph:\n# w:hanging=\"244\" (-12.2pt) -> w:hanging=\"210\" (-10.5pt)\nforeach ($p in $d.Paragraphs) {\n    $p.Range.ParagraphFormat.FirstLineIndent = -10.5\n}\n\n# Drop the trailing space in the first bullet's text:\n# \"Unordered information. \" -> \"Unordered information.\"\n$find = $d.Content.Find\n$find.Execute(\"Unordered information. \", $false, $false, $false, $false, $false, $true, 1, $false, \"Unordered information.\", 2)\n"}
